# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to the freshly scraped counts (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$exhibitionUpdates = @{
    3  = 541
    4  = 1556
    5  = 159
    8  = 173
    9  = 754
    10 = 1052
    11 = 68
    12 = 340
    14 = 3
    15 = 14
    16 = 6490
    22 = 15498
    23 = 1537
    24 = 290
    25 = 148
    27 = 11101
    28 = 765
    29 = 4346
    30 = 253
    33 = 307
    34 = 127
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new F value
$allTypesUpdates = @{
    3  = 541
    4  = 1556
    5  = 159
    9  = 173
    10 = 754
    12 = 1052
    13 = 68
    14 = 340
    16 = 3
    18 = 14
    19 = 6490
    26 = 15498
    27 = 1537
    28 = 290
    29 = 148
    32 = 11101
    33 = 765
    34 = 4346
    35 = 253
    38 = 307
    39 = 127
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
